$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Row 2 (National Trends, 2025, C_0) updates
$ws.Range("J2").Value = 68
$ws.Range("K2").Value = -600
$ws.Range("O2").Value = 100
$ws.Range("P2").Value = 68

# Row 4 (Distributed Energy, 2040, C_0) updates
$ws.Range("O4").Value = 329
$ws.Range("P4").Value = 214
$ws.Range("W4").Value = 1890

# Row 6 (Distributed Energy, 2030, C_0) updates
$ws.Range("H6").Value = 216
$ws.Range("I6").Value = 172
$ws.Range("J6").Value = 128
$ws.Range("L6").Value = 970
$ws.Range("O6").Value = 197
$ws.Range("P6").Value = 130
$ws.Range("R6").Value = 437
$ws.Range("S6").Value = 411
$ws.Range("T6").Value = 298
$ws.Range("U6").Value = 225
$ws.Range("V6").Value = 194
$ws.Range("X6").Value = 183
$ws.Range("Y6").Value = 207
$ws.Range("Z6").Value = 900
$ws.Range("AA6").Value = 240

# Update the active cell selection to Z13
$ws.Activate()
$ws.Range("Z13").Select()
